$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5436
$ws.Range("E2").Value = 119
$ws.Range("F2").Value = 119
$ws.Range("G2").Value = 86
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 3633
$ws.Range("L2").Value = 2162
$ws.Range("M2").Value = 1472
$ws.Range("N2").Value = 1451
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 298
$ws.Range("Q2").Value = 101
$ws.Range("R2").Value = -149
$ws.Range("S2").Value = 18
$ws.Range("T2").Value = 78
$ws.Range("U2").Value = 23
$ws.Range("V2").Value = 1265
$ws.Range("W2").Value = 2.19
$ws.Range("X2").Value = 1.22
$ws.Range("Y2").Value = 4.4
$ws.Range("Z2").Value = 1.85
$ws.Range("AA2").Value = 146.91
$ws.Range("AB2").Value = 379.89
$ws.Range("AC2").Value = 1048
$ws.Range("AD2").Value = 11.35
$ws.Range("AE2").Value = 24329
$ws.Range("AF2").Value = 0.49
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 1.68
$ws.Range("AI2").Value = 19.1
$ws.Range("AJ2").Value = 5960575

# Row 3
$ws.Range("D3").Value = 7151
$ws.Range("E3").Value = 178
$ws.Range("F3").Value = 178
$ws.Range("G3").Value = 194
$ws.Range("H3").Value = 141
$ws.Range("I3").Value = 122
$ws.Range("J3").Value = 18
$ws.Range("K3").Value = 6570
$ws.Range("L3").Value = 4078
$ws.Range("M3").Value = 2493
$ws.Range("N3").Value = 1757
$ws.Range("O3").Value = 735
$ws.Range("P3").Value = 298
$ws.Range("Q3").Value = 183
$ws.Range("R3").Value = 19
$ws.Range("S3").Value = -163
$ws.Range("T3").Value = 151
$ws.Range("U3").Value = 32
$ws.Range("V3").Value = 2529
$ws.Range("W3").Value = 2.49
$ws.Range("X3").Value = 1.97
$ws.Range("Y3").Value = 7.63
$ws.Range("Z3").Value = 2.76
$ws.Range("AA3").Value = 163.6
$ws.Range("AB3").Value = 482.67
$ws.Range("AC3").Value = 2052
$ws.Range("AD3").Value = 8.97
$ws.Range("AE3").Value = 29465
$ws.Range("AF3").Value = 0.62
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 5960575

# Row 4
$ws.Range("D4").Value = 9536
$ws.Range("E4").Value = 364
$ws.Range("F4").Value = 364
$ws.Range("G4").Value = 217
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 68
$ws.Range("J4").Value = 32
$ws.Range("K4").Value = 6717
$ws.Range("L4").Value = 4126
$ws.Range("M4").Value = 2592
$ws.Range("N4").Value = 1829
$ws.Range("O4").Value = 763
$ws.Range("P4").Value = 298
$ws.Range("Q4").Value = 605
$ws.Range("R4").Value = -373
$ws.Range("S4").Value = -183
$ws.Range("T4").Value = 161
$ws.Range("U4").Value = 445
$ws.Range("V4").Value = 2343
$ws.Range("W4").Value = 3.82
$ws.Range("X4").Value = 1.05
$ws.Range("Y4").Value = 3.82
$ws.Range("Z4").Value = 1.51
$ws.Range("AA4").Value = 159.2
$ws.Range("AB4").Value = 505.6
$ws.Range("AC4").Value = 1147
$ws.Range("AD4").Value = 15.78
$ws.Range("AE4").Value = 30659
$ws.Range("AF4").Value = 0.59
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 0.55
$ws.Range("AI4").Value = 8.72
$ws.Range("AJ4").Value = 5960575

# Row 5
$ws.Range("D5").Value = 9874
$ws.Range("E5").Value = 554
$ws.Range("F5").Value = 554
$ws.Range("G5").Value = 574
$ws.Range("H5").Value = 423
$ws.Range("I5").Value = 287
$ws.Range("J5").Value = 136
$ws.Range("K5").Value = 6801
$ws.Range("L5").Value = 3783
$ws.Range("M5").Value = 3018
$ws.Range("N5").Value = 2123
$ws.Range("O5").Value = 895
$ws.Range("P5").Value = 298
$ws.Range("Q5").Value = 518
$ws.Range("R5").Value = -235
$ws.Range("S5").Value = -279
$ws.Range("T5").Value = 209
$ws.Range("U5").Value = 310
$ws.Range("V5").Value = 2030
$ws.Range("W5").Value = 5.61
$ws.Range("X5").Value = 4.29
$ws.Range("Y5").Value = 14.53
$ws.Range("Z5").Value = 6.26
$ws.Range("AA5").Value = 125.37
$ws.Range("AB5").Value = 607.29
$ws.Range("AC5").Value = 4813
$ws.Range("AD5").Value = 5.11
$ws.Range("AE5").Value = 35596
$ws.Range("AF5").Value = 0.69
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 0.61
$ws.Range("AI5").Value = 3.12
$ws.Range("AJ5").Value = 5960575

# Row 6
$ws.Range("D6").Value = 9371
$ws.Range("E6").Value = 363
$ws.Range("F6").Value = 363
$ws.Range("G6").Value = 203
$ws.Range("H6").Value = 144
$ws.Range("I6").Value = 99
$ws.Range("K6").Value = 7209
$ws.Range("L6").Value = 4046
$ws.Range("M6").Value = 3163
$ws.Range("N6").Value = 2226
$ws.Range("P6").Value = 298
$ws.Range("Q6").Value = 76
$ws.Range("R6").Value = -313
$ws.Range("S6").Value = 281
$ws.Range("T6").Value = 362
$ws.Range("U6").Value = -286
$ws.Range("V6").Value = 2323
$ws.Range("W6").Value = 3.87
$ws.Range("X6").Value = 1.53
$ws.Range("Y6").Value = 4.56
$ws.Range("Z6").Value = 2.05
$ws.Range("AA6").Value = 127.9
$ws.Range("AB6").Value = 638.88
$ws.Range("AC6").Value = 1663
$ws.Range("AD6").Value = 14.55
$ws.Range("AE6").Value = 37327
$ws.Range("AF6").Value = 0.65
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 0.62
$ws.Range("AI6").Value = 9.02
$ws.Range("AJ6").Value = 5960575

# Row 7: clear all data columns D:AJ
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all data columns D:AJ
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all data columns D:AJ
$ws.Range("D9:AJ9").ClearContents()
